$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("G1").Value = "variants"
$ws.Range("H1").Value = "GATK-filtered SNPs"

# --- Remove the mostly-empty "STACKS -- myco run -- SE" placeholder row (old row 8) ---
$ws.Rows(8).Delete()

# Relabel the data-bearing STACKS row (old row 7, still row 7 after the delete) from PE to SE
$ws.Range("B7").Value = "STACKS -- myco run -- SE"

# At this point:
#   row 8 = "STACKS -- pub reported" (was row 9)
#   row 9 = "ISSRseq" full data row (was row 10)
# Swap them (via a scratch row far below the data) so the ISSRseq row ends up at 8
# and the pub-reported row ends up at 9, matching the target layout.
$ws.Range("A9:J9").Cut($ws.Range("A20:J20"))
$ws.Range("A8:J8").Cut($ws.Range("A9:J9"))
$ws.Range("A20:J20").Cut($ws.Range("A8:J8"))

# Fix up number formats so they reflect each row's actual data styling:
# the ISSRseq row (now row 8) uses thousands-separator/center format on G:I,
$ws.Range("G8:I8").NumberFormat = "#,##0"
# but plain center format on J.
$ws.Range("J2").Copy()
$ws.Range("J8").PasteSpecial(-4122)
# the "STACKS -- pub reported" row (now row 9) uses plain center format on G:H.
$ws.Range("G2:H2").Copy()
$ws.Range("G9:H9").PasteSpecial(-4122)
$ws.Range("I9:J9").Clear()

# Clear the scratch row entirely.
$ws.Range("A20:J20").Clear()

# Selection, as saved by Excel after editing
$ws.Range("G1").Select()
